# Edit script for NIT-9005959406.xlsx
# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" table (rows 16-112, 97 rows ascending 1706..2506) is
# replaced with a table listing periods in descending order (2507..1706),
# i.e. the newest period 2507 is added and the whole list is reversed so
# the most recent period is now on top. The table grows to 98 rows
# (16-113). The last row keeps the special "total row" formatting/border
# that used to belong to row 112, and the two signature footer rows move
# down from 117/118 to 118/119.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Grow the data table by one row -------------------------------
# Copy the old last row (112, special border/style) down to the new last
# row (113) before we touch anything else, so we keep its formatting.
$ws.Range("B112:J112").Copy($ws.Range("B113:J113"))

# Revert row 112 back to the regular (non-last-row) look by copying the
# formatting/values from a normal data row (16).
$ws.Range("B16:J16").Copy($ws.Range("B112:J112"))

# --- 2. Move the two footer/signature rows down by one ---------------
$ws.Range("B118:C118").Copy($ws.Range("B119:C119"))
$ws.Range("H118:J118").Copy($ws.Range("H119:J119"))
$ws.Range("B117:C117").Copy($ws.Range("B118:C118"))
$ws.Range("H117:J117").Copy($ws.Range("H118:J118"))
$ws.Range("B117:C117").UnMerge()
$ws.Range("H117:J117").UnMerge()
$ws.Range("B117:C117").Clear()
$ws.Range("H117:J117").Clear()

# --- 3. Re-fill the "Periodo Mora" column (E16:E113) with the new,    -
#        descending period list, newest (2507) first, oldest (1706)    -
#        last.                                                         -
$periods = @("2507","2506","2505","2504","2503","2502","2501","2412","2411","2410","2409","2408","2407","2406","2405","2404","2403","2402","2401","2312","2311","2310","2309","2308","2307","2306","2305","2304","2303","2302","2301","2212","2211","2210","2209","2208","2207","2206","2205","2204","2203","2202","2201","2112","2111","2110","2109","2108","2107","2106","2105","2104","2103","2102","2101","2012","2011","2010","2009","2008","2007","2006","2005","2004","2003","2002","2001","1912","1911","1910","1909","1908","1907","1906","1905","1904","1903","1902","1901","1812","1811","1810","1809","1808","1807","1806","1805","1804","1803","1802","1801","1712","1711","1710","1709","1708","1707","1706")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# --- 4. "Valor Mora" (column F): every row is 41459 except the very   -
#        last one (now row 113, the oldest period 1706) which keeps     -
#        the original 29021.                                           -
for ($row = 16; $row -le 112; $row++) {
    $ws.Range("F$row").Value = 41459
}
$ws.Range("F113").Value = 29021

# --- 5. Header / summary cells ----------------------------------------
$ws.Range("E11").Value = 4050544   # VALOR MORA total
$ws.Range("F13").Value = 98        # Cant. Periodos
